$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.990.20'
$ws.Range("E2").Value = '  -1.64%  '
$ws.Range("D3").Value = '1.907.84'
$ws.Range("E3").Value = '  -3.15%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.93'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.70%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.13%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4599'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -1.43%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3826'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -2.30%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07718'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -2.73%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9805'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.78%  '
$ws.Range("E11").Value = '  -3.20%  '
$ws.Range("D12").Value = '1.907.72'
$ws.Range("E12").Value = '  -4.50%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.668'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -2.30%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.933'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -3.47%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07064'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.71%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.005'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.15%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '83.73'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -4.59%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009470'
$ws.Range("D18").ClearFormats()
$ws.Range("E19").Value = '  -3.31%  '
$ws.Range("E20").Value = '  -0.23%  '
$ws.Range("D21").Value = '28.956.87'
$ws.Range("E21").Value = '  -1.99%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.314'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -3.94%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.88'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -2.36%  '
$ws.Range("D24").Value = '2.141.55'
$ws.Range("E24").Value = '  -4.33%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.094'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.83%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '158.28'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.19%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.09'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.15%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.647'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -2.60%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '117.44'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.96%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.852'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -2.39%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09298'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.34%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.8615'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -3.13%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.079'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -2.98%  '
$ws.Range("E34").Value = '  -5.63%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.975'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -6.36%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05730'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.57%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.151'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -1.50%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.003'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.01%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02040'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -3.17%  '
$ws.Range("E40").Value = '  -3.57%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.401'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -4.43%  '
$ws.Range("E42").Value = '  -2.22%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.846'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +3.63%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '9.318'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -3.37%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5184'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -2.65%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '11.22'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -4.61%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.06826'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.44%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.053'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -5.01%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '111.07'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -2.02%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.777'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -2.94%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.000002479'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -14.36%  '
